$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 389; Date = "2021-09-04"; CumCases = 6941611; NewCases = 37578; NewDeaths = 120; CumDeaths = 133161 },
    @{ Row = 390; Date = "2021-09-05"; CumCases = 6978126; NewCases = 37011; NewDeaths = 68;  CumDeaths = 133229 },
    @{ Row = 391; Date = "2021-09-06"; CumCases = 7018927; NewCases = 41192; NewDeaths = 45;  CumDeaths = 133274 },
    @{ Row = 392; Date = "2021-09-07"; CumCases = 7056106; NewCases = 37489; NewDeaths = 209; CumDeaths = 133483 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A holds an ISO date string as literal text (matches existing rows),
    # so force Text format before assignment to stop Excel auto-converting it
    # to a date serial number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.Date

    $ws.Cells.Item($row, 2).Value = "overview"
    $ws.Cells.Item($row, 3).Value = "K02000001"
    $ws.Cells.Item($row, 4).Value = "United Kingdom"
    $ws.Cells.Item($row, 5).Value = $r.CumCases
    $ws.Cells.Item($row, 6).Value = $r.NewCases
    $ws.Cells.Item($row, 7).Value = $r.NewDeaths
    $ws.Cells.Item($row, 8).Value = $r.CumDeaths
}
